# Update the "Handback" report timestamps to reflect the newly generated
# report (commit message: "Generate Report for Handback").
#
# Mapping of changed shared-string values (old -> new):
#   Overview!G2  2016-08-30 01:06:30 -> 2016-08-30 01:07:21
#   zh-cn!H2     2016-08-30 01:06:26 -> 2016-08-30 01:07:16
#   zh-cn!K2     2016-08-30 01:06:50 -> 2016-08-30 01:07:33
#   de-de!K2     2016-08-30 01:06:57 -> 2016-08-30 01:07:41

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-30 01:07:21"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-30 01:07:16"
$wsZhCn.Range("K2").Value = "2016-08-30 01:07:33"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-30 01:07:41"
